# The sheet had a header row (date / error_id / description / station_id)
# in row 1 followed by two data rows. To make the data export cleanly to
# postgres (per the commit message), the header row is removed so the
# table starts directly with the data, shifting what were rows 2 and 3 up
# to rows 1 and 2.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the header row; Excel shifts the rows below it up automatically.
$ws.Rows("1:1").Delete() | Out-Null

# Leave the UI selection on the now-first row, matching the post-edit
# saved state (a whole-row selection, as after a row deletion).
$ws.Rows("1:1").Select() | Out-Null
